$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: "===================Criar nova DB========================"
#         -> "===================Criar nova DB 1========================"
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "===================Criar nova DB========================",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "===================Criar nova DB 1========================", 2)

# ---------------------------------------------------------------------------
# Step 2: remove the _GoBack bookmark that currently sits inside the
#         "Computador(...)" paragraph -- it will be re-created as its own
#         paragraph at the very end of the document.
# ---------------------------------------------------------------------------
$bk = $d.Bookmarks.Item("_GoBack")
$bk.Delete()

# ---------------------------------------------------------------------------
# Step 3: append the new "Criar nova DB 2" block just before the very last
#         (trailing) empty paragraph of the document.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($n - 1).Range

$anchor.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$p1 = $d.Paragraphs.Item($n - 1)
$p1.Range.Text = "===================Criar nova DB 2========================"

$anchor = $d.Paragraphs.Item($n - 1).Range
$anchor.InsertParagraphAfter()
$n = $d.Paragraphs.Count

$anchor = $d.Paragraphs.Item($n - 1).Range
$anchor.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$p3 = $d.Paragraphs.Item($n - 1)
$p3.Range.Text = "Enderecos(id, rua, cidade, cep)"
$t3 = $p3.Range.Text
$idx3 = $t3.IndexOf("id")
$s3 = $p3.Range.Start + $idx3
$r3 = $d.Range($s3, $s3 + 2)
$r3.Font.Bold = $true

$anchor = $d.Paragraphs.Item($n - 1).Range
$anchor.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$p4 = $d.Paragraphs.Item($n - 1)
$p4.Range.Text = "Site(id, nome, idEndereco)"
$t4 = $p4.Range.Text
$idx4 = $t4.IndexOf("id, ")
$s4 = $p4.Range.Start + $idx4
$r4 = $d.Range($s4, $s4 + 4)
$r4.Font.Bold = $true

$anchor = $d.Paragraphs.Item($n - 1).Range
$anchor.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$p5 = $d.Paragraphs.Item($n - 1)
$d.Bookmarks.Add("_GoBack", $p5.Range)

Write-Output "done"
